{"js": "// Use case 2 - Indstil lydpakker: rework list/term-explanation wording.\n// 1. \"Sekund\u00e6r akt\u00f8r\" value: Webserver -> Ingen\n// 2. \"Pr\u00e6kondition\" value: En USB- eller netv\u00e6rksforbindelse er tilkoblet -> Der er oprettet forbindelse til USB-drev\n// 3. Remove the \"Bruger v\u00e6lger USB eller server som kilde\" bullet from Hovedscenarie\n// 4. Renumber \"Undtagelse 4.a\" -> \"Undtagelse 3.a\" (both occurrences) after the removal above\n// 5. Reword the first bullet under \"Undtagelse 2.a\" (lydpakkeliste wording)\n\nconst body = context.document.body;\n\n// 1) Sekund\u00e6r akt\u00f8r: Webserver -> Ingen\nconst secondaryActor = body.search(\"Webserver\", { matchCase: true });\nsecondaryActor.load(\"text\");\nawait context.sync();\nfor (const item of secondaryActor.items) {\n  item.insertText(\"Ingen\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Pr\u00e6kondition value\nconst precondition = body.search(\"En USB- eller netv\u00e6rksforbindelse er tilkoblet\", { matchCase: true });\nprecondition.load(\"text\");\nawait context.sync();\nfor (const item of precondition.items) {\n  item.insertText(\"Der er oprettet forbindelse til USB-drev\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 3) Remove the whole bullet paragraph \"Bruger v\u00e6lger USB eller server som kilde\"\nconst sourceBullet = body.search(\"Bruger v\u00e6lger USB eller server som kilde\", { matchCase: true });\nsourceBullet.load(\"text\");\nawait context.sync();\nfor (const item of sourceBullet.items) {\n  const para = item.paragraphs.getFirst();\n  para.delete();\n}\nawait context.sync();\n\n// 4) Renumber both \"Undtagelse 4.a\" references to \"Undtagelse 3.a\"\nconst exceptionRef = body.search(\" 4.a: Der \", { matchCase: true });\nexceptionRef.load(\"text\");\nawait context.sync();\nfor (const item of exceptionRef.items) {\n  item.insertText(\" 3.a: Der \", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 5) Reword the \"lydpakkeliste\" bullet to explain the list inline\nconst listWording = body.search(\"Bruger v\u00e6lger lydpakke p\u00e5 lydpakkeliste\", { matchCase: true });\nlistWording.load(\"text\");\nawait context.sync();\nfor (const item of listWording.items) {\n  item.insertText(\"Bruger v\u00e6lger lydpakke fra liste over systemets lydpakker\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Use case 2 - Indstil lydpakker: rework list/term-explanation wording.\n# 1. \"Sekund\u00e6r akt\u00f8r\" value: Webserver -> Ingen\n# 2. \"Pr\u00e6kondition\" value: En USB- eller netv\u00e6rksforbindelse er tilkoblet -> Der er oprettet forbindelse til USB-drev\n# 3. Remove the \"Bruger v\u00e6lger USB eller server som kilde\" bullet from Hovedscenarie\n# 4. Renumber \"Undtagelse 4.a\" -> \"Undtagelse 3.a\" (both occurrences) after the removal above\n# 5. Reword the first bullet under \"Undtagelse 2.a\" (lydpakkeliste wording)\n\n$d = $word.ActiveDocument\n\n# 1) Sekund\u00e6r akt\u00f8r: Webserver -> Ingen\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Webserver\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ingen\", 2) | Out-Null\n\n# 2) Pr\u00e6kondition value\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"En USB- eller netv\u00e6rksforbindelse er tilkoblet\", $false, $false, $false, $false, $false, $true, 1, $false, \"Der er oprettet forbindelse til USB-drev\", 2) | Out-Null\n\n# 3) Remove the whole bullet paragraph \"Bruger v\u00e6lger USB eller server som kilde\"\n$range = $d.Content\n$found = $range.Find.Execute(\"Bruger v\u00e6lger USB eller server som kilde\")\nif ($found) {\n    $para = $range.Paragraphs(1)\n    $para.Range.Delete()\n}\n\n# 4) Renumber both \"Undtagelse 4.a\" references to \"Undtagelse 3.a\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\" 4.a: Der \", $false, $false, $false, $false, $false, $true, 1, $false, \" 3.a: Der \", 2) | Out-Null\n\n# 5) Reword the \"lydpakkeliste\" bullet to explain the list inline\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Bruger v\u00e6lger lydpakke p\u00e5 lydpakkeliste\", $false, $false, $false, $false, $false, $true, 1, $false, \"Bruger v\u00e6lger lydpakke fra liste over systemets lydpakker\", 2) | Out-Null\n"}
